$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.031580648947838
$ws.Cells.Item(2, 4).Value = 1.035612323162507
$ws.Cells.Item(2, 5).Value = 1.035219840682011
$ws.Cells.Item(2, 6).Value = 1.030184597301502
$ws.Cells.Item(2, 9).Value = 1.035353984675241
$ws.Cells.Item(2, 10).Value = 1.03671558986735
$ws.Cells.Item(2, 11).Value = 1.038408428788788
$ws.Cells.Item(2, 12).Value = 1.0380170713037
$ws.Cells.Item(2, 13).Value = 1.032996343376906
$ws.Cells.Item(2, 14).Value = 1.038187844917921

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.032584800311452
$ws.Cells.Item(3, 4).Value = 1.036374725853028
$ws.Cells.Item(3, 5).Value = 1.036172015830103
$ws.Cells.Item(3, 6).Value = 1.03182024049123
$ws.Cells.Item(3, 9).Value = 1.035604357769741
$ws.Cells.Item(3, 10).Value = 1.037361397781079
$ws.Cells.Item(3, 11).Value = 1.038980387374563
$ws.Cells.Item(3, 12).Value = 1.038778216293469
$ws.Cells.Item(3, 13).Value = 1.034438066176892
$ws.Cells.Item(3, 14).Value = 1.038834569952944

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.033234373658651
$ws.Cells.Item(4, 4).Value = 1.036867847721962
$ws.Cells.Item(4, 5).Value = 1.036788308065755
$ws.Cells.Item(4, 6).Value = 1.032878376671704
$ws.Cells.Item(4, 9).Value = 1.035765088307554
$ws.Cells.Item(4, 10).Value = 1.037778519231655
$ws.Cells.Item(4, 11).Value = 1.039349619701264
$ws.Cells.Item(4, 12).Value = 1.039270280846643
$ws.Cells.Item(4, 13).Value = 1.035370261352346
$ws.Cells.Item(4, 14).Value = 1.039252283763823

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.033507412127436
$ws.Cells.Item(5, 4).Value = 1.037075107434477
$ws.Cells.Item(5, 5).Value = 1.037047438738481
$ws.Cells.Item(5, 6).Value = 1.033323167423243
$ws.Cells.Item(5, 9).Value = 1.035832353730503
$ws.Cells.Item(5, 10).Value = 1.037953695828553
$ws.Cells.Item(5, 11).Value = 1.039504638641469
$ws.Cells.Item(5, 12).Value = 1.039477038237693
$ws.Cells.Item(5, 13).Value = 1.035761995141987
$ws.Cells.Item(5, 14).Value = 1.039427709131582

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.033553254040098
$ws.Cells.Item(6, 4).Value = 1.037109904396018
$ws.Cells.Item(6, 5).Value = 1.037090950359491
$ws.Cells.Item(6, 6).Value = 1.033397847019023
$ws.Cells.Item(6, 9).Value = 1.035843629969628
$ws.Cells.Item(6, 10).Value = 1.037983098136771
$ws.Cells.Item(6, 11).Value = 1.039530654913103
$ws.Cells.Item(6, 12).Value = 1.039511747479951
$ws.Cells.Item(6, 13).Value = 1.035827759684792
$ws.Cells.Item(6, 14).Value = 1.039457153194451

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.033238022176647
$ws.Cells.Item(7, 4).Value = 1.036870617327798
$ws.Cells.Item(7, 5).Value = 1.036791770419232
$ws.Cells.Item(7, 6).Value = 1.032884320174821
$ws.Cells.Item(7, 9).Value = 1.035765988312398
$ws.Cells.Item(7, 10).Value = 1.037780860661466
$ws.Cells.Item(7, 11).Value = 1.039351691882297
$ws.Cells.Item(7, 12).Value = 1.039273043967063
$ws.Cells.Item(7, 13).Value = 1.035375496344409
$ws.Cells.Item(7, 14).Value = 1.039254628518734

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.031920044336975
$ws.Cells.Item(8, 4).Value = 1.035870022911402
$ws.Cells.Item(8, 5).Value = 1.035541597529167
$ws.Cells.Item(8, 6).Value = 1.03073742295325
$ws.Cells.Item(8, 9).Value = 1.035438864068815
$ws.Cells.Item(8, 10).Value = 1.036934001260248
$ws.Cells.Item(8, 11).Value = 1.038601903770697
$ws.Cells.Item(8, 12).Value = 1.038274396706177
$ws.Cells.Item(8, 13).Value = 1.033483727767978
$ws.Cells.Item(8, 14).Value = 1.038406566480029

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.029596183274567
$ws.Cells.Item(9, 4).Value = 1.034105285617163
$ws.Cells.Item(9, 5).Value = 1.033339925084832
$ws.Cells.Item(9, 6).Value = 1.026952234288878
$ws.Cells.Item(9, 9).Value = 1.034852636832519
$ws.Cells.Item(9, 10).Value = 1.035435884910696
$ws.Cells.Item(9, 11).Value = 1.037274056610166
$ws.Cells.Item(9, 12).Value = 1.03651120485671
$ws.Cells.Item(9, 13).Value = 1.030144612309823
$ws.Cells.Item(9, 14).Value = 1.036906322633425

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.028045933049078
$ws.Cells.Item(10, 4).Value = 1.032927737448151
$ws.Cells.Item(10, 5).Value = 1.031872981529407
$ws.Cells.Item(10, 6).Value = 1.02442697395725
$ws.Cells.Item(10, 9).Value = 1.034455222681977
$ws.Cells.Item(10, 10).Value = 1.034433175137988
$ws.Cells.Item(10, 11).Value = 1.036384348886981
$ws.Cells.Item(10, 12).Value = 1.035333390549597
$ws.Cells.Item(10, 13).Value = 1.027914437661363
$ws.Cells.Item(10, 14).Value = 1.035902188897828

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.027374404409513
$ws.Cells.Item(11, 4).Value = 1.032417593090921
$ws.Cells.Item(11, 5).Value = 1.031237969789554
$ws.Cells.Item(11, 6).Value = 1.02333299146078
$ws.Cells.Item(11, 9).Value = 1.034281570369356
$ws.Cells.Item(11, 10).Value = 1.033998041029314
$ws.Cells.Item(11, 11).Value = 1.035998028229941
$ws.Cells.Item(11, 12).Value = 1.034822817225199
$ws.Cells.Item(11, 13).Value = 1.026947697810813
$ws.Cells.Item(11, 14).Value = 1.035466436848809

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.027124928412889
$ws.Cells.Item(12, 4).Value = 1.032228063657148
$ws.Cells.Item(12, 5).Value = 1.031002125100668
$ws.Cells.Item(12, 6).Value = 1.022926550530923
$ws.Cells.Item(12, 9).Value = 1.034216832058781
$ws.Cells.Item(12, 10).Value = 1.033836268619458
$ws.Cells.Item(12, 11).Value = 1.035854370001414
$ws.Cells.Item(12, 12).Value = 1.034633080791507
$ws.Cells.Item(12, 13).Value = 1.026588441650498
$ws.Cells.Item(12, 14).Value = 1.035304434703575

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.027178443708193
$ws.Cells.Item(13, 4).Value = 1.032268720138733
$ws.Cells.Item(13, 5).Value = 1.031052713380233
$ws.Cells.Item(13, 6).Value = 1.023013737557771
$ws.Cells.Item(13, 9).Value = 1.034230729335968
$ws.Cells.Item(13, 10).Value = 1.03387097589127
$ws.Cells.Item(13, 11).Value = 1.035885192498288
$ws.Cells.Item(13, 12).Value = 1.034673783840576
$ws.Cells.Item(13, 13).Value = 1.026665510934162
$ws.Cells.Item(13, 14).Value = 1.035339191263695

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.027353783476623
$ws.Cells.Item(14, 4).Value = 1.032401927328362
$ws.Cells.Item(14, 5).Value = 1.031218474248835
$ws.Cells.Item(14, 6).Value = 1.023299396749096
$ws.Cells.Item(14, 9).Value = 1.0342762238989
$ws.Cells.Item(14, 10).Value = 1.033984671818742
$ws.Cells.Item(14, 11).Value = 1.035986156694606
$ws.Cells.Item(14, 12).Value = 1.034807135327313
$ws.Cells.Item(14, 13).Value = 1.02691800500876
$ws.Cells.Item(14, 14).Value = 1.035453048652424

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.027461810708107
$ws.Cells.Item(15, 4).Value = 1.032483995483032
$ws.Cells.Item(15, 5).Value = 1.031320608537965
$ws.Cells.Item(15, 6).Value = 1.023475389027106
$ws.Cells.Item(15, 9).Value = 1.034304223303851
$ws.Cells.Item(15, 10).Value = 1.034054704494213
$ws.Cells.Item(15, 11).Value = 1.036048342645732
$ws.Cells.Item(15, 12).Value = 1.034889286056673
$ws.Cells.Item(15, 13).Value = 1.02707355275074
$ws.Cells.Item(15, 14).Value = 1.035523180782328

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.028090494604243
$ws.Cells.Item(16, 4).Value = 1.032961588582249
$ws.Cells.Item(16, 5).Value = 1.031915128988597
$ws.Cells.Item(16, 6).Value = 1.024499566182822
$ws.Cells.Item(16, 9).Value = 1.034466714321217
$ws.Cells.Item(16, 10).Value = 1.034462033378076
$ws.Cells.Item(16, 11).Value = 1.036409965095517
$ws.Cells.Item(16, 12).Value = 1.03536726349783
$ws.Cells.Item(16, 13).Value = 1.02797857410272
$ws.Cells.Item(16, 14).Value = 1.035931088119926

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.028484780870041
$ws.Cells.Item(17, 4).Value = 1.033261100765238
$ws.Cells.Item(17, 5).Value = 1.032288104813343
$ws.Cells.Item(17, 6).Value = 1.025141858586368
$ws.Cells.Item(17, 9).Value = 1.034568220248032
$ws.Cells.Item(17, 10).Value = 1.034717283954323
$ws.Cells.Item(17, 11).Value = 1.036636514069917
$ws.Cells.Item(17, 12).Value = 1.03566693240134
$ws.Cells.Item(17, 13).Value = 1.028545981667249
$ws.Cells.Item(17, 14).Value = 1.036186701181269

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.028714736372813
$ws.Cells.Item(18, 4).Value = 1.033435776260136
$ws.Cells.Item(18, 5).Value = 1.032505673155713
$ws.Cells.Item(18, 6).Value = 1.025516446447464
$ws.Cells.Item(18, 9).Value = 1.034627275539383
$ws.Cells.Item(18, 10).Value = 1.034866075293517
$ws.Cells.Item(18, 11).Value = 1.036768552924063
$ws.Cells.Item(18, 12).Value = 1.035841669073883
$ws.Cells.Item(18, 13).Value = 1.028876839127223
$ws.Cells.Item(18, 14).Value = 1.036335703821232

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.028793141054417
$ws.Cells.Item(19, 4).Value = 1.033495331876662
$ws.Cells.Item(19, 5).Value = 1.032579861389732
$ws.Cells.Item(19, 6).Value = 1.025644162909975
$ws.Cells.Item(19, 9).Value = 1.034647386197828
$ws.Cells.Item(19, 10).Value = 1.034916793716208
$ws.Cells.Item(19, 11).Value = 1.036813557252111
$ws.Cells.Item(19, 12).Value = 1.035901240438151
$ws.Cells.Item(19, 13).Value = 1.028989635913618
$ws.Cells.Item(19, 14).Value = 1.036386494269901

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.028442480280695
$ws.Cells.Item(20, 4).Value = 1.033228968519544
$ws.Cells.Item(20, 5).Value = 1.032248086195814
$ws.Cells.Item(20, 6).Value = 1.025072951961747
$ws.Cells.Item(20, 9).Value = 1.034557345284879
$ws.Cells.Item(20, 10).Value = 1.034689907505649
$ws.Cells.Item(20, 11).Value = 1.036612218200904
$ws.Cells.Item(20, 12).Value = 1.035634786472129
$ws.Cells.Item(20, 13).Value = 1.028485114780349
$ws.Cells.Item(20, 14).Value = 1.036159285854898

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.027302151409575
$ws.Cells.Item(21, 4).Value = 1.032362702213149
$ws.Cells.Item(21, 5).Value = 1.031169661066957
$ws.Cells.Item(21, 6).Value = 1.023215279704739
$ws.Cells.Item(21, 9).Value = 1.034262833404495
$ws.Cells.Item(21, 10).Value = 1.033951195190358
$ws.Cells.Item(21, 11).Value = 1.035956429719113
$ws.Cells.Item(21, 12).Value = 1.034767869040153
$ws.Cells.Item(21, 13).Value = 1.026843656432092
$ws.Cells.Item(21, 14).Value = 1.035419524483388

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.026584947233764
$ws.Cells.Item(22, 4).Value = 1.031817819792721
$ws.Cells.Item(22, 5).Value = 1.03049176789212
$ws.Cells.Item(22, 6).Value = 1.022046779491933
$ws.Cells.Item(22, 9).Value = 1.034076295871784
$ws.Cells.Item(22, 10).Value = 1.033485902542698
$ws.Cells.Item(22, 11).Value = 1.035543174585516
$ws.Cells.Item(22, 12).Value = 1.034222301210299
$ws.Cells.Item(22, 13).Value = 1.025810642998118
$ws.Cells.Item(22, 14).Value = 1.034953571066798

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.026965173365186
$ws.Cells.Item(23, 4).Value = 1.032106693840467
$ws.Cells.Item(23, 5).Value = 1.030851117411219
$ws.Cells.Item(23, 6).Value = 1.022666274317716
$ws.Cells.Item(23, 9).Value = 1.034175312560852
$ws.Cells.Item(23, 10).Value = 1.033732642412025
$ws.Cells.Item(23, 11).Value = 1.035762337762646
$ws.Cells.Item(23, 12).Value = 1.034511564946367
$ws.Cells.Item(23, 13).Value = 1.026358356550976
$ws.Cells.Item(23, 14).Value = 1.035200661335042

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.028461594173077
$ws.Cells.Item(24, 4).Value = 1.03324348777498
$ws.Cells.Item(24, 5).Value = 1.032266168831033
$ws.Cells.Item(24, 6).Value = 1.025104088053397
$ws.Cells.Item(24, 9).Value = 1.034562259680282
$ws.Cells.Item(24, 10).Value = 1.034702278028381
$ws.Cells.Item(24, 11).Value = 1.03662319677735
$ws.Cells.Item(24, 12).Value = 1.035649312003877
$ws.Cells.Item(24, 13).Value = 1.028512618219436
$ws.Cells.Item(24, 14).Value = 1.036171673945191

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.030197130742664
$ws.Cells.Item(25, 4).Value = 1.034561698245728
$ws.Cells.Item(25, 5).Value = 1.033908960197357
$ws.Cells.Item(25, 6).Value = 1.027931084223375
$ws.Cells.Item(25, 9).Value = 1.035005352320136
$ws.Cells.Item(25, 10).Value = 1.035823880053503
$ws.Cells.Item(25, 11).Value = 1.0376181242518
$ws.Cells.Item(25, 12).Value = 1.036967444365062
$ws.Cells.Item(25, 13).Value = 1.031008550371833
$ws.Cells.Item(25, 14).Value = 1.037294868773838
